$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("B4", 0.305),
    @("D4", 0.223),
    @("E4", 0.194),
    @("F4", 0.027),
    @("G4", 0.166),
    @("H4", 0.206),
    @("I4", 0.027),
    @("J4", 0.166),
    @("K4", 0.346),
    @("L4", 0.1),
    @("M4", 0.317),
    @("N4", 0.255),
    @("O4", 0.023),
    @("P4", 0.153),
    @("Q4", 0.5),
    @("R4", 0.231),
    @("T4", 0.235),
    @("U4", 0.08799999999999999),
    @("V4", 0.297),
    @("W4", 0.228),
    @("Z4", 0.437),
    @("AA4", 0.135),
    @("AB4", 0.368),
    @("AC4", 0.118),
    @("AE4", 0.082),
    @("AF4", 0.714),
    @("AH4", 0.328),
    @("AI4", 0.667),
    @("AJ4", 0.158),
    @("AK4", 0.398),
    @("AL4", 0.646),
    @("AO4", 0.676),
    @("B5", 0.821),
    @("C5", 0.147),
    @("D5", 0.384),
    @("E5", 0.769),
    @("F5", 0.178),
    @("G5", 0.421),
    @("H5", 0.821),
    @("I5", 0.147),
    @("J5", 0.384),
    @("K5", 0.667),
    @("L5", 0.222),
    @("M5", 0.471),
    @("N5", 0.795),
    @("O5", 0.163),
    @("P5", 0.404),
    @("Q5", 0.538),
    @("R5", 0.249),
    @("S5", 0.499),
    @("T5", 0.487),
    @("W5", 0.6919999999999999),
    @("X5", 0.213),
    @("Y5", 0.462),
    @("Z5", 0.795),
    @("AA5", 0.163),
    @("AB5", 0.404),
    @("AC5", 0.718),
    @("AD5", 0.202),
    @("AE5", 0.45),
    @("AF5", 0.949),
    @("AG5", 0.049),
    @("AH5", 0.221),
    @("AI5", 0.795),
    @("AJ5", 0.163),
    @("AK5", 0.404),
    @("AL5", 0.897),
    @("AM5", 0.092),
    @("AN5", 0.303),
    @("AO5", 0.88),
    @("B6", 0.445),
    @("E6", 0.31),
    @("H6", 0.329),
    @("K6", 0.456),
    @("N6", 0.386),
    @("Q6", 0.518),
    @("T6", 0.317),
    @("W6", 0.343),
    @("Z6", 0.5639999999999999),
    @("AC6", 0.203),
    @("AF6", 0.8149999999999999),
    @("AI6", 0.725),
    @("AL6", 0.751),
    @("AO6", 0.764),
    @("B7", 0.613),
    @("E7", 0.483),
    @("H7", 0.514),
    @("K7", 0.5629999999999999),
    @("N7", 0.5580000000000001),
    @("Q7", 0.53),
    @("T7", 0.401),
    @("W7", 0.492),
    @("Z7", 0.6830000000000001),
    @("AC7", 0.356),
    @("AF7", 0.89),
    @("AI7", 0.766),
    @("AL7", 0.832),
    @("AO7", 0.829),
    @("B8", 0.769),
    @("C8", 0.153),
    @("D8", 0.391),
    @("E8", 0.649),
    @("F8", 0.166),
    @("G8", 0.407),
    @("H8", 0.72),
    @("I8", 0.152),
    @("J8", 0.39),
    @("K8", 0.587),
    @("L8", 0.199),
    @("M8", 0.446),
    @("N8", 0.72),
    @("O8", 0.159),
    @("P8", 0.399),
    @("Q8", 0.529),
    @("R8", 0.243),
    @("S8", 0.493),
    @("T8", 0.443),
    @("V8", 0.471),
    @("W8", 0.629),
    @("X8", 0.197),
    @("Y8", 0.444),
    @("Z8", 0.735),
    @("AA8", 0.159),
    @("AB8", 0.398),
    @("AC8", 0.594),
    @("AD8", 0.187),
    @("AE8", 0.433),
    @("AF8", 0.889),
    @("AG8", 0.063),
    @("AH8", 0.25),
    @("AI8", 0.795),
    @("AJ8", 0.163),
    @("AK8", 0.404),
    @("AL8", 0.86),
    @("AM8", 0.097),
    @("AN8", 0.311),
    @("AO8", 0.848),
    @("B9", 0.718),
    @("C9", 0.202),
    @("D9", 0.45),
    @("E9", 0.513),
    @("H9", 0.615),
    @("I9", 0.237),
    @("J9", 0.487),
    @("K9", 0.487),
    @("L9", 0.25),
    @("M9", 0.5),
    @("N9", 0.615),
    @("O9", 0.237),
    @("P9", 0.487),
    @("Q9", 0.513),
    @("R9", 0.25),
    @("S9", 0.5),
    @("T9", 0.385),
    @("U9", 0.237),
    @("V9", 0.487),
    @("W9", 0.538),
    @("X9", 0.249),
    @("Y9", 0.499),
    @("Z9", 0.641),
    @("AA9", 0.23),
    @("AB9", 0.48),
    @("AC9", 0.487),
    @("AD9", 0.25),
    @("AE9", 0.5),
    @("AF9", 0.795),
    @("AG9", 0.163),
    @("AH9", 0.404),
    @("AI9", 0.795),
    @("AJ9", 0.163),
    @("AK9", 0.404),
    @("AL9", 0.795),
    @("AM9", 0.163),
    @("AN9", 0.404),
    @("AO9", 0.795),
    @("B10", 0.769),
    @("C10", 0.178),
    @("D10", 0.421),
    @("E10", 0.6919999999999999),
    @("F10", 0.213),
    @("G10", 0.462),
    @("H10", 0.744),
    @("I10", 0.191),
    @("J10", 0.437),
    @("K10", 0.667),
    @("L10", 0.222),
    @("M10", 0.471),
    @("N10", 0.769),
    @("O10", 0.178),
    @("P10", 0.421),
    @("Q10", 0.538),
    @("R10", 0.249),
    @("S10", 0.499),
    @("T10", 0.487),
    @("W10", 0.6919999999999999),
    @("X10", 0.213),
    @("Y10", 0.462),
    @("Z10", 0.795),
    @("AA10", 0.163),
    @("AB10", 0.404),
    @("AC10", 0.59),
    @("AD10", 0.242),
    @("AE10", 0.492),
    @("AF10", 0.949),
    @("AG10", 0.049),
    @("AH10", 0.221),
    @("AI10", 0.795),
    @("AJ10", 0.163),
    @("AK10", 0.404),
    @("AL10", 0.897),
    @("AM10", 0.092),
    @("AN10", 0.303),
    @("AO10", 0.88),
    @("B11", 0.821),
    @("C11", 0.147),
    @("D11", 0.384),
    @("E11", 0.769),
    @("F11", 0.178),
    @("G11", 0.421),
    @("H11", 0.821),
    @("I11", 0.147),
    @("J11", 0.384),
    @("K11", 0.667),
    @("L11", 0.222),
    @("M11", 0.471),
    @("N11", 0.795),
    @("O11", 0.163),
    @("P11", 0.404),
    @("Q11", 0.538),
    @("R11", 0.249),
    @("S11", 0.499),
    @("T11", 0.487),
    @("W11", 0.6919999999999999),
    @("X11", 0.213),
    @("Y11", 0.462),
    @("Z11", 0.795),
    @("AA11", 0.163),
    @("AB11", 0.404),
    @("AC11", 0.641),
    @("AD11", 0.23),
    @("AE11", 0.48),
    @("AF11", 0.949),
    @("AG11", 0.049),
    @("AH11", 0.221),
    @("AI11", 0.795),
    @("AJ11", 0.163),
    @("AK11", 0.404),
    @("AL11", 0.897),
    @("AM11", 0.092),
    @("AN11", 0.303),
    @("AO11", 0.88),
    @("B12", 1.281),
    @("C12", 0.64),
    @("D12", 0.8),
    @("E12", 1.633),
    @("F12", 1.032),
    @("G12", 1.016),
    @("H12", 1.562),
    @("I12", 1.309),
    @("J12", 1.144),
    @("K12", 1.423),
    @("L12", 0.552),
    @("M12", 0.743),
    @("N12", 1.323),
    @("O12", 0.477),
    @("P12", 0.6899999999999999),
    @("Z12", 1.226),
    @("AA12", 0.239),
    @("AB12", 0.489),
    @("AC12", 2.107),
    @("AD12", 4.239),
    @("AE12", 2.059),
    @("AF12", 1.189),
    @("AG12", 0.207),
    @("AH12", 0.455),
    @("AL12", 1.114),
    @("AM12", 0.101),
    @("AN12", 0.318),
    @("AO12", 1.101),
    @("B13", 3.385),
    @("C13", 1.416),
    @("D13", 1.19),
    @("E13", 4.564),
    @("F13", 0.707),
    @("G13", 0.841),
    @("H13", 4.526),
    @("I13", 0.9340000000000001),
    @("J13", 0.966),
    @("K13", 2.324),
    @("L13", 0.572),
    @("M13", 0.756),
    @("N13", 3.308),
    @("O13", 0.828),
    @("P13", 0.91),
    @("Z13", 2.838),
    @("AA13", 4.19),
    @("AB13", 2.047),
    @("AC13", 6.184),
    @("AD13", 2.94),
    @("AE13", 1.715),
    @("AF13", 1.641),
    @("AG13", 0.743),
    @("AH13", 0.862),
    @("AI13", 1.333),
    @("AJ13", 0.376),
    @("AK13", 0.613),
    @("AL13", 1.744),
    @("AM13", 0.857),
    @("AN13", 0.926),
    @("AO13", 1.573)
)

foreach ($chg in $changes) {
    $ws.Range($chg[0]).Value = $chg[1]
}

Write-Host "Applied" $changes.Count "cell updates"